# Apply "moved local logs to repo" edits to the Activity Log sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header section ---
$ws.Range("B2").Value = "Jesse Hare"
$ws.Range("F2").Value = "Week 3"

# --- Activity descriptions (column A), entered top to bottom ---
$ws.Range("A4").Value = "Meeting with Client"
$ws.Range("A5").Value = "Elicict and analyse requirements"
$ws.Range("A6").Value = "Finalise project Plan"
$ws.Range("A7").Value = "Finalise project Plan"
$ws.Range("A8").Value = "Finalise project Plan"
$ws.Range("A9").Value = "Finalise project Plan"

# --- Type (column C), entered top to bottom ---
$ws.Range("C4").Value = "G"
$ws.Range("C5").Value = "G"
$ws.Range("C6").Value = "G"
$ws.Range("C7").Value = "G"
$ws.Range("C8").Value = "G"
$ws.Range("C9").Value = "G"

# --- Dates (column D) ---
$ws.Range("D4").Value = 43689
$ws.Range("D5").Value = 43689
$ws.Range("D6").Value = 43690
$ws.Range("D7").Value = 43691
$ws.Range("D8").Value = 43692
$ws.Range("D9").Value = 43693

# --- Start times (column E) ---
$ws.Range("E4").Value = 0.375
$ws.Range("E5").Value = 0.41666666666666669
$ws.Range("E6").Value = 0.375
$ws.Range("E7").Value = 0.375
$ws.Range("E8").Value = 0.375
$ws.Range("E9").Value = 0.375

# --- End times (column F) ---
$ws.Range("F4").Value = 0.41666666666666669
$ws.Range("F5").Value = 0.45833333333333331
$ws.Range("F6").Value = 0.58333333333333337
$ws.Range("F7").Value = 0.54166666666666663
$ws.Range("F8").Value = 0.083333333333333329
$ws.Range("F9").Value = 0.54166666666666663

# --- Hours (column G) ---
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 1
$ws.Range("G6").Value = 5
$ws.Range("G7").Value = 4
$ws.Range("G8").Value = 5
$ws.Range("G9").Value = 4

# --- Column widths, to mirror Excel's auto-sizing of the data columns ---
$ws.Columns.Item(2).ColumnWidth = 22.15
$ws.Columns.Item(4).ColumnWidth = 12.15
$ws.Columns.Item(5).ColumnWidth = 12.83
$ws.Columns.Item(6).ColumnWidth = 12.83

# Selection / active cell like the saved file
$ws.Range("F7").Select()

$wb.Save()
